# Adds a "Col4" date column (D) to the worksheet, populating rows 2-21
# with sequential dates formatted as short dates, matching the
# "Added a date support for excel parser" change (issue #791).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new column
$ws.Range("D1").Value = "Col4"

# Sequential date serials (2014-05-20 .. 2014-06-08), one per data row
$dates = 41779,41780,41781,41782,41783,41784,41785,41786,41787,41788,41789,41790,41791,41792,41793,41794,41795,41796,41797,41798

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 4).Value = $dates[$i]
}

# Apply the date format to D2 then replicate the same style to the
# remaining date cells (copy/paste-format keeps a single shared xf,
# same as re-using one style record for the whole column).
$ws.Range("D2").NumberFormat = "mm-dd-yy"
$ws.Range("D2").Copy()
$ws.Range("D3:D21").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Resize column D to fit its new contents
$ws.Columns("D").AutoFit()

# Reflect the new selection/active cell like the authored workbook
$ws.Range("D1:D21").Select()

# Keep the printed page as A4, matching the resaved workbook
$ws.PageSetup.PaperSize = 9
